$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3 through 17 (the bulk of old data), keeping header row 1 and data row 2
$ws.Range("A3:A17").EntireRow.Delete()

# Update the remaining data row with the new values
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 6.681483765882756
